$d = $word.ActiveDocument

# Step 1: Remove the trailing run's text "IOT/Temp_Monitor/Code" first,
# while it is still unique in the document (before the merge below would
# otherwise create a duplicate substring inside the combined URL).
$d.Content.Find.Execute(
    "IOT/Temp_Monitor/Code",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2
)

# Step 2: Replace the first URL run's text with the combined full URL.
$d.Content.Find.Execute(
    "https://github.com/Devan2120/Networking/",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://github.com/Devan2120/Networking/tree/main/IOT/Temp_Monitor/Code",
    2
)
